$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '54.501.46'
$ws.Range("E2").Value = '  -5.66%  '

# Row 3
$ws.Range("D3").Value = '2.877.01'
$ws.Range("E3").Value = '  -9.10%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '469.04'
$ws.Range("E5").Value = '  -11.40%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '125.59'
$ws.Range("E6").Value = '  -6.89%  '

# Row 7
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("D8").Value = '2.874.95'
$ws.Range("E8").Value = '  -9.16%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.406'
$ws.Range("E9").Value = '  -10.28%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.61'
$ws.Range("E10").Value = '  -8.11%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0971'
$ws.Range("E11").Value = '  -13.21%  '

# Row 12
$ws.Range("E12").Value = '  -16.11%  '

# Row 13
$ws.Range("E13").Value = '  -4.23%  '

# Row 14
$ws.Range("D14").Value = '3.368.08'
$ws.Range("E14").Value = '  -9.19%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.35'
$ws.Range("E15").Value = '  -9.94%  '

# Row 16
$ws.Range("D16").Value = '54.429.24'
$ws.Range("E16").Value = '  -5.69%  '

# Row 17
$ws.Range("D17").Value = '2.872.41'
$ws.Range("E17").Value = '  -9.31%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000133'
$ws.Range("E18").Value = '  -13.47%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.35'
$ws.Range("E19").Value = '  -8.60%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.43'
$ws.Range("E20").Value = '  -13.41%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.06'
$ws.Range("E21").Value = '  -12.88%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '300.14'
$ws.Range("E22").Value = '  -14.15%  '

# Row 23
$ws.Range("E23").Value = '  -0.20%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.439'
$ws.Range("E24").Value = '  -14.66%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '58.85'
$ws.Range("E25").Value = '  -15.64%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  -0.26%  '

# Row 27
$ws.Range("E27").Value = '  -10.32%  '

# Row 28
$ws.Range("E28").Value = '  +0.05%  '

# Row 29
$ws.Range("D29").Value = '0.0₃0798'
$ws.Range("E29").Value = '  -17.56%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.12'
$ws.Range("E30").Value = '  -12.01%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.20'
$ws.Range("E31").Value = '  -11.19%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.10'
$ws.Range("E32").Value = '  -9.79%  '

# Row 33
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.81'
$ws.Range("E33").Value = '  -13.36%  '

# Row 34
$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.60'
$ws.Range("E34").Value = '  -15.43%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '141.56'
$ws.Range("E35").Value = '  -11.13%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.17'
$ws.Range("E36").Value = '  -16.43%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.39'
$ws.Range("E37").Value = '  -14.21%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.21'
$ws.Range("E38").Value = '  -15.02%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.85'
$ws.Range("E39").Value = '  -13.23%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0619'
$ws.Range("E40").Value = '  -11.94%  '

# Row 41
$ws.Range("D41").Value = '2.906.50'
$ws.Range("E41").Value = '  -9.03%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.23%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '35.08'
$ws.Range("E43").Value = '  -13.07%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.603'
$ws.Range("E44").Value = '  -13.71%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.945'
$ws.Range("E45").Value = '  -14.04%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.42'
$ws.Range("E46").Value = '  -14.37%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.31'
$ws.Range("E47").Value = '  -11.58%  '

# Row 48
$ws.Range("D48").Value = '2.029.43'
$ws.Range("E48").Value = '  -10.84%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.34'
$ws.Range("E49").Value = '  -14.39%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0214'
$ws.Range("E50").Value = '  -9.09%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.70'
$ws.Range("E51").Value = '  -14.85%  '
